$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value - process different datatypes (isinstance-style change)
$ws.Range("B3").Value = 9.99999

# Update the active selection to B3
$ws.Range("B3").Select()
